# Revert "Ajout de légendes et d'une option pour quitter"
#
# The journal table (Tableau1) previously grew from E5:M28 to E5:M40:
#   - row 28 was filled in with a "Légendes" entry
#   - row 29 was filled in with a "quitter" entry
#   - rows 30-40 were appended as fresh blank placeholder rows
# This reverts that: drop rows 29-40 entirely, and restore row 28 back to
# being an empty placeholder row (same shape as the other blank rows that
# already existed below the real data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 29 through 40 — the table (Tableau1) and the sheet dimension
# auto-shrink to E5:M28 as part of this.
$ws.Rows("29:40").Delete()

# Row 28 (formerly the "Légendes" entry) goes back to being blank, like the
# placeholder rows used to be. Keep the H28 formula (it's part of every row
# in the table) but clear the data cells around it.
$ws.Range("E28").ClearContents()
$ws.Range("F28").ClearContents()
$ws.Range("G28").ClearContents()
$ws.Range("I28").ClearContents()
$ws.Range("J28").ClearContents()
$ws.Range("K28").ClearContents()
$ws.Range("L28").ClearContents()
$ws.Range("M28").ClearContents()

# Restore the plain (non-date/time) number format + centered alignment that
# the other blank rows use for their Date/Heure Début/Heure fin cells, by
# copying the formatting already used on the (also blank) I28 cell.
$ws.Range("I28").Copy()
$ws.Range("E28:G28").PasteSpecial(-4122)
$ws.Range("K28").PasteSpecial(-4122)

# Row 28 no longer has wrapped text, so let it size back down to the
# default row height instead of staying at the old wrapped-text height.
$ws.Rows("28").AutoFit()

# Restore the selection to where it ended up after the revert.
$ws.Range("L28").Select()
